$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header column: author (F)
$ws.Range("F1").Value = "author"

# Fill in author for each data row, in row order so "jah" is interned
# before "creation_date"
$ws.Range("F2").Value = "jah"
$ws.Range("F3").Value = "jah"
$ws.Range("F4").Value = "jah"

# New header column: creation_date (G)
$ws.Range("G1").Value = "creation_date"

# Apply the date number format to G2 first (built-in "mm-dd-yy", numFmtId
# 14) then copy that formatting onto G3:G4 so all three cells share the
# exact same style record instead of each getting its own duplicate.
$ws.Range("G2").NumberFormat = "mm-dd-yy"
$ws.Range("G2").Copy()
$ws.Range("G3:G4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$creationDate = Get-Date -Year 2022 -Month 3 -Day 11 -Hour 0 -Minute 0 -Second 0
$ws.Range("G2").Value = $creationDate
$ws.Range("G3").Value = $creationDate
$ws.Range("G4").Value = $creationDate

# Size the new column to fit its ("creation_date") header.
$ws.Columns.Item(7).ColumnWidth = 10.83

# Re-fit the wrapped rows now that the row content/format has changed.
$ws.Rows.Item(2).RowHeight = 48
$ws.Rows.Item(3).RowHeight = 128
$ws.Rows.Item(4).RowHeight = 48

$ws.Range("G3").Select()
